$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp shown in the header row.
$ws.Range("F1").Value = "Last status check on: 24.02.2022 11:45"

# Row 9 (OMV IKEA) price refresh:
#  - current price (B9) bumped to 38.9
#  - previous price (C9) now holds the prior current price (38.5)
#  - delta (D9) becomes a plain text label "+0.4" (was a numeric 0.6)
#  - the "old date" (E9) becomes a literal text timestamp (was a numeric
#    Excel date/time with a date-time number format) and loses that
#    number formatting, reverting to the sheet's default (unstyled) cell
$ws.Range("B9").Value = 38.9
$ws.Range("C9").Value = 38.5

$plainStyle = $ws.Range("A2").Style

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "+0.4"
$ws.Range("D9").Style = $plainStyle

$ws.Range("E9").Value = "2022-02-24 11:48:56"
$ws.Range("E9").Style = $plainStyle
